# Apply updated Betfair Back/Lay odds values for 2025-12-11 workbook
# Generated from the authoritative cell-level diff between before/after OOXML
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("Q2").Value = 1.99

# Row 3
$ws.Range("F3").Value = 3.6
$ws.Range("W3").Value = 1.36
$ws.Range("Y3").Value = 10.5
$ws.Range("AI3").Value = 34
$ws.Range("AO3").Value = 15.5

# Row 4
$ws.Range("F4").Value = 4.8
$ws.Range("G4").Value = 5.1
$ws.Range("H4").Value = 1.84
$ws.Range("I4").Value = 1.87
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 3.65
$ws.Range("P4").Value = 1.92
$ws.Range("R4").Value = 1.35
$ws.Range("S4").Value = 3.65
$ws.Range("T4").Value = 1.92
$ws.Range("V4").Value = 2.14
$ws.Range("X4").Value = 13
$ws.Range("Y4").Value = 8.4
$ws.Range("AA4").Value = 19.5
$ws.Range("AC4").Value = 8.199999999999999
$ws.Range("AF4").Value = 38
$ws.Range("AJ4").Value = 130
$ws.Range("AK4").Value = 70

# Row 5
$ws.Range("K5").Value = 4.1
$ws.Range("Q5").Value = 1.82
$ws.Range("S5").Value = 3.05
$ws.Range("T5").Value = 1.78
$ws.Range("AJ5").Value = 20
$ws.Range("AO5").Value = 50

# Row 6
$ws.Range("N6").Value = 5.4
$ws.Range("Q6").Value = 1.63
$ws.Range("S6").Value = 2.54

# Row 7
$ws.Range("V7").Value = 1.83
$ws.Range("AG7").Value = 15
$ws.Range("AO7").Value = 15.5

# Row 8
$ws.Range("G8").Value = 3.2
$ws.Range("H8").Value = 2.34
$ws.Range("I8").Value = 2.42
$ws.Range("O8").Value = 1.19
$ws.Range("P8").Value = 2.58
$ws.Range("Q8").Value = 1.59
$ws.Range("R8").Value = 1.62
$ws.Range("S8").Value = 2.46
$ws.Range("T8").Value = 1.53
$ws.Range("U8").Value = 2.72
$ws.Range("V8").Value = 1.7
$ws.Range("AC8").Value = 9.4
$ws.Range("AG8").Value = 14
$ws.Range("AH8").Value = 14.5
$ws.Range("AN8").Value = 20
$ws.Range("AO8").Value = 12.5

# Row 9
$ws.Range("J9").Value = 7.4
$ws.Range("T9").Value = 1.97
$ws.Range("U9").Value = 1.95

# Row 10
$ws.Range("N10").Value = 5

# Row 11
$ws.Range("F11").Value = 3.7
$ws.Range("H11").Value = 2.14
$ws.Range("S11").Value = 3.6
$ws.Range("T11").Value = 1.83

# Row 12
$ws.Range("J12").Value = 3.2
$ws.Range("K12").Value = 3.7
$ws.Range("V12").Value = 1.65

# Row 13
$ws.Range("F13").Value = 8
$ws.Range("G13").Value = 9.199999999999999
$ws.Range("H13").Value = 1.42
$ws.Range("I13").Value = 1.44
$ws.Range("K13").Value = 5.6
$ws.Range("R13").Value = 1.5
$ws.Range("S13").Value = 2.74
$ws.Range("V13").Value = 3.3
$ws.Range("W13").Value = 1.12
$ws.Range("Z13").Value = 8.800000000000001

# Row 14
$ws.Range("Q14").Value = 1.86

# Row 15
$ws.Range("X15").Value = 16

# Row 16
$ws.Range("G16").Value = 2.04
$ws.Range("K16").Value = 4
$ws.Range("W16").Value = 1.97

# Row 17
$ws.Range("AJ17").Value = 55

# Row 18
$ws.Range("F18").Value = 1.48
$ws.Range("G18").Value = 1.52
$ws.Range("M18").Value = 1.05
$ws.Range("Q18").Value = 1.8
$ws.Range("S18").Value = 2.74
$ws.Range("X18").Value = 22
$ws.Range("AK18").Value = 15

# Row 19
$ws.Range("H19").Value = 2.44
$ws.Range("V19").Value = 1.58
$ws.Range("AG19").Value = 14.5

# Row 20
$ws.Range("T20").Value = 1.85
$ws.Range("U20").Value = 2
$ws.Range("AF20").Value = 10

# Row 21
$ws.Range("G21").Value = 5.4
$ws.Range("N21").Value = 4.8
$ws.Range("P21").Value = 2.3
$ws.Range("Q21").Value = 1.68
$ws.Range("S21").Value = 2.7
$ws.Range("T21").Value = 1.74

# Row 22
$ws.Range("J22").Value = 7.6
$ws.Range("P22").Value = 2.64
$ws.Range("Q22").Value = 1.53
$ws.Range("R22").Value = 1.66
$ws.Range("T22").Value = 2.24
$ws.Range("AD22").Value = 1000
$ws.Range("AN22").Value = 3.95

# Row 23
$ws.Range("F23").Value = 2.44
$ws.Range("I23").Value = 3.3
$ws.Range("L23").Value = 1.44
$ws.Range("N23").Value = 3.55
$ws.Range("R23").Value = 1.33
$ws.Range("W23").Value = 1.64

# Row 24
$ws.Range("F24").Value = 4.5
$ws.Range("H24").Value = 1.79
$ws.Range("P24").Value = 2.52
$ws.Range("Q24").Value = 1.59
$ws.Range("R24").Value = 1.62

# Row 25
$ws.Range("G25").Value = 1.25
$ws.Range("H25").Value = 14
$ws.Range("J25").Value = 7.4
$ws.Range("K25").Value = 8.800000000000001
$ws.Range("N25").Value = 4.7
$ws.Range("U25").Value = 1.58
$ws.Range("V25").Value = 1.06
$ws.Range("W25").Value = 5
$ws.Range("Z25").Value = 200
$ws.Range("AD25").Value = 70
$ws.Range("AE25").Value = 420
$ws.Range("AI25").Value = 290

# Row 26
$ws.Range("F26").Value = 3.3
$ws.Range("G26").Value = 3.4
$ws.Range("H26").Value = 2.36
$ws.Range("I26").Value = 2.4
$ws.Range("J26").Value = 3.55
$ws.Range("Q26").Value = 1.92
$ws.Range("S26").Value = 3.3
$ws.Range("V26").Value = 1.72
$ws.Range("W26").Value = 1.42
$ws.Range("X26").Value = 13.5
$ws.Range("AA26").Value = 32
$ws.Range("AG26").Value = 13.5
$ws.Range("AL26").Value = 46
$ws.Range("AN26").Value = 34
$ws.Range("AO26").Value = 18.5

# Row 28
$ws.Range("S28").Value = 2.66

# Row 29
$ws.Range("F29").Value = 22
$ws.Range("G29").Value = 29
$ws.Range("J29").Value = 8
$ws.Range("K29").Value = 9.800000000000001
$ws.Range("N29").Value = 5.7
$ws.Range("P29").Value = 2.62
$ws.Range("T29").Value = 2.42
$ws.Range("U29").Value = 1.57
$ws.Range("AA29").Value = 8.4
$ws.Range("AG29").Value = 1000
$ws.Range("AI29").Value = 1000
$ws.Range("AK29").Value = 1000

# Row 30
$ws.Range("AD30").Value = 15

# Row 31
$ws.Range("I31").Value = 2.6
$ws.Range("T31").Value = 1.68

# Row 32
$ws.Range("P32").Value = 2.26
$ws.Range("AN32").Value = 29

# Row 33
$ws.Range("P33").Value = 1.94
$ws.Range("AB33").Value = 10.5
$ws.Range("AD33").Value = 16

# Row 34
$ws.Range("S34").Value = 2.96
$ws.Range("AF34").Value = 55
$ws.Range("AI34").Value = 40
$ws.Range("AK34").Value = 1000

# Row 35
$ws.Range("H35").Value = 2.2
$ws.Range("I35").Value = 2.26
$ws.Range("J35").Value = 3.7
$ws.Range("L35").Value = 1.32
$ws.Range("N35").Value = 4.5
$ws.Range("S35").Value = 2.9
$ws.Range("V35").Value = 1.79
$ws.Range("Z35").Value = 15.5
$ws.Range("AA35").Value = 48
$ws.Range("AC35").Value = 8.6
$ws.Range("AE35").Value = 23
$ws.Range("AF35").Value = 27
$ws.Range("AG35").Value = 15.5
$ws.Range("AJ35").Value = 1000
$ws.Range("AK35").Value = 1000
$ws.Range("AN35").Value = 30

# Row 37
$ws.Range("N37").Value = 3.7
$ws.Range("U37").Value = 2.06
$ws.Range("AK37").Value = 1000
